$d = $word.ActiveDocument

# Change 1: insert "在作者的心目中，" after "可以初步窥探"
$d.Content.Find.Execute(
    "可以初步窥探徽宗是怎样",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "可以初步窥探在作者的心目中，徽宗是怎样",
    2
)

# Change 2: rewrite "徽宗大部分都在参与和政事无太大关系的活动，" -> "徽宗的大部分活动，都是在参与和政事无太大关系的事业，"
$d.Content.Find.Execute(
    "徽宗大部分都在参与和政事无太大关系的活动，",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "徽宗的大部分活动，都是在参与和政事无太大关系的事业，",
    2
)

# Change 3: "对于评价徽宗，关键性的问题，便是北宋覆灭的责任归属。" -> "评价徽宗，关键性的问题，是北宋覆灭的责任归属。"
$d.Content.Find.Execute(
    "对于评价徽宗，关键性的问题，便是北宋覆灭的责任归属。",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "评价徽宗，关键性的问题，是北宋覆灭的责任归属。",
    2
)

# Change 4: touch the last paragraph so its paragraph-mark run properties
# drop the eastAsia hint (mirrors what Word does when the cursor visits
# the end of that paragraph during editing).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Font.Name = $lastPara.Range.Font.Name
